# feat: add fifteenth lesson
#
# Adds two new rows to the hour log (lessons 15 and 16), mirroring the
# existing table rows 3-16 (Lesson No. / Date / Hours columns A:C).
#   Row 17 -> Lesson 15, 2017-12-18 (serial 43087), 2 hours
#   Row 18 -> Lesson 16, 2017-12-19 (serial 43088), 2 hours
# F2 (Remaining hours = F1 - SUM(C:C)) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: fifteenth lesson -------------------------------------------------
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 43087
$ws.Cells.Item(17, 3).Value = 2

# --- Row 18: sixteenth lesson --------------------------------------------------
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 43088
$ws.Cells.Item(18, 3).Value = 2

# Carry over the same formatting used by the rest of the table (date format
# on column B, centered alignment on column C) from the last existing row.
$ws.Range("A16:C16").Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16:C16").Copy() | Out-Null
$ws.Range("A18:C18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# PasteSpecial(formats) should not touch values, but re-assert them to be safe.
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 43087
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 43088
$ws.Cells.Item(18, 3).Value = 2

# Move the active selection onto the newly-entered dates, as a user would
# leave it right after typing the new rows.
$ws.Range("B17:B18").Select() | Out-Null
